{"js": "// Amendments from final proof reading:\n// In the RQ1 \"Answer\" table cell, the closing sentence about maintainability\n// is reworded from:\n//   \"The impact on maintainability is inconclusive with indicators towards\n//    both enhanced and degraded maintainability.\"\n// to:\n//   \"These observations are consistent with degradation in maintainability.\"\n// The rest of the paragraph (\"Greater team sizes result in degraded\n// measures ... enhanced measures modularity. \") is left untouched.\n\nconst oldTail =\n  \"The impact on maintainability is inconclusive with indicators towards both enhanced and degraded maintainability.\";\nconst newTail =\n  \"These observations are consistent with degradation in maintainability.\";\n\nconst body = context.document.body;\nconst results = body.search(oldTail, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found: \" + oldTail);\n}\n\n// Replace the matched range's text in place; the surrounding run formatting\n// (color/size) is preserved because insertText(\"Replace\") reuses the\n// formatting of the range being replaced.\nresults.items[0].insertText(newTail, \"Replace\");\nawait context.sync();\n", "ps1": "# Amendments from final proof reading\n#\n# In the RQ1 \"Answer\" table cell, reword the closing sentence about\n# maintainability from:\n#   \"The impact on maintainability is inconclusive with indicators towards\n#    both enhanced and degraded maintainability.\"\n# to:\n#   \"These observations are consistent with degradation in maintainability.\"\n#\n# The rest of the paragraph (\"Greater team sizes result in degraded\n# measures ... enhanced measures modularity. \") is left untouched, and the\n# run formatting (color/size) carries over from the text being replaced.\n\n$d = $word.ActiveDocument\n\n$oldTail = \"The impact on maintainability is inconclusive with indicators towards both enhanced and degraded maintainability.\"\n$newTail = \"These observations are consistent with degradation in maintainability.\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute($oldTail, $false, $true, $false, $false, $false, $true, 1, $false, $newTail, 2)\n\nif (-not $found) {\n    throw \"Target sentence not found: $oldTail\"\n}\n"}
